$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Insert a new row at position 8 to hold the newly scraped listing,
# shifting the previous rows 8-16 down to 9-17.
$ws.Rows("8:8").Insert()

# Refresh the acquisition timestamp for every data row (2-17).
$newTimestamp = '2025-10-28 18:28:07'
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Populate row 8 with the new listing data (columns B-E, G, H; F/hyperlink handled below).
$ws.Range("B8").Value = '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("G8").Value = 135
$ws.Range("H8").Value = '◆ツール,スクレイピング ◇サイト'

# Rebuild the hyperlinks on column F: the row insert does not shift the
# existing hyperlink-to-cell associations, so clear them all out and
# re-add one per data row, in order, pointing at the correct URLs.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5416301') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5416307') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5416305') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5422389') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5421873') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5416328') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5251319') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5421687') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5421820') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5411585') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5421779') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5422004') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5422125') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), 'https://www.lancers.jp/work/detail/5421982') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), 'https://www.lancers.jp/work/detail/5421894') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), 'https://www.lancers.jp/work/detail/5422200') | Out-Null

# Re-use the original hyperlink cell style (as seen on a still-clean cell)
# instead of the fresh duplicate style that Hyperlinks.Add() just created,
# so every URL cell keeps a single, shared "Hyperlink" style.
$linkStyle = $ws.Cells.Item(2, 6).Style
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).Style = $linkStyle
}

# Widen column H slightly (17 -> 19) to match the new content.
# (ColumnWidth uses Excel character units, which are offset from the raw
#  OOXML column width by 5/6 in this engine, hence the 18.1667 below.)
$ws.Columns("H").ColumnWidth = 18.166666666666668

$ws.Range("A1").Select()
